$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.882.15'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '2.935.28'
$ws.Range("E3").Value = '  +4.09%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.93'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.11'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.562'
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.41'
$ws.Range("E10").Value = '  -1.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").Value = '  +5.30%  '
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.00'
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.77'
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.394.28'
$ws.Range("E15").Value = '  +3.82%  '
$ws.Range("D16").Value = '2.937.97'
$ws.Range("E16").Value = '  +4.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.986'
$ws.Range("E17").Value = '  +1.45%  '
$ws.Range("D18").Value = '51.899.55'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.63'
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("B20").Value = 'ImmutableX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.32'
$ws.Range("E20").Value = '  -3.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.22'
$ws.Range("E21").Value = '  +6.64%  '
$ws.Range("D22").Value = '0.0₃0988'
$ws.Range("E22").Value = '  +1.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.26'
$ws.Range("E23").Value = '  +1.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.71'
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.79'
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("E26").Value = '  +10.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.93'
$ws.Range("E27").Value = '  +3.03%  '
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.32'
$ws.Range("E29").Value = '  +16.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.106'
$ws.Range("E30").Value = '  +18.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.58'
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.37'
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.27'
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.21'
$ws.Range("E34").Value = '  +11.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '52.86'
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0454'
$ws.Range("E36").Value = '  +1.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.32'
$ws.Range("E38").Value = '  +3.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.89'
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.05'
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("E41").Value = '  +7.46%  '
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.06'
$ws.Range("E43").Value = '  +4.64%  '
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("E45").Value = '  +1.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.52'
$ws.Range("E46").Value = '  +1.22%  '
$ws.Range("D47").Value = '2.172.94'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '111.64'
$ws.Range("E48").Value = '  -9.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.248'
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0347'
$ws.Range("E50").Value = '  +12.72%  '
$ws.Range("E51").Value = '  -0.66%  '
